$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: N10 updated, O10 updated
$ws.Range("N10").Value = 177670.86
$ws.Range("O10").Value = 177555.36

# Row 15: O15 filled in
$ws.Range("O15").Value = 1684.05

# Row 21: N21 updated
$ws.Range("N21").Value = 720590.33

# Row 29: N29 updated
$ws.Range("N29").Value = 202098

# Row 30: N30 updated
$ws.Range("N30").Value = 18940.7

# Row 34: N34 updated
$ws.Range("N34").Value = 31962.1
